# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Price (col D) and Volume(1h) (col E) are kept as plain text, matching the
# sheet's existing convention (values like "2.398.02" aren't valid numbers,
# and the % values carry intentional padding spaces). For cells whose new
# text would otherwise auto-coerce to a number (e.g. "0.998"), the range is
# briefly formatted as Text so the literal string sticks, then restored to
# the default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.955.07"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.398.02"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.427.78"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0971"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.27%  "
$ws.Range("D14").Value = "2.817.45"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "56.733.97"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "2.370.21"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "2.488.89"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.380"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.150"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").Value = "0.0₃0736"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "
# Stacks and OKB swapped ranking positions (rows 42/43).
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "258.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.567"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0494"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0212"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
